$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet was protected (with only the data-entry column B cells left
# unlocked). The author unprotected it as part of this edit.
$ws.Unprotect()

# 2. Информация об организации / Организация: the reporting organization
# changed from the Ministry of Digital Development to the State Agency of
# Communications.
$ws.Range("B6").Value = "Государственное агентство связи при Государственном комитете информационных технологий и связи Кыргызской Республики"

# A new (empty) formatted cell was introduced at C3, extending the used
# range from A1:B26 to A1:C26.
$ws.Range("C3").Style = "Обычный 2"
$ws.Range("C3").Font.Name = "Calibri"

# Final UI selection moved from B3 to B4.
$ws.Range("B4").Select()
